# "we just use actual bools"
# Add a new "boolean" row to the Main sheet's settings table, using real
# boolean cells (t="b") instead of numeric 0/1 placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "boolean"
$ws.Range("B29").Value = $true
$ws.Range("C29").Value = $true
$ws.Range("D29").Value = $false

# Reflect where the editor's selection ended up after adding the row.
$null = $ws.Range("C39").Select()
